$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8-10 (old MuSCs sending-cluster rows no longer present after TPM recalculation)
$ws.Rows("8:10").Delete()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "F10"
$ws.Range("C2").Value = "F3"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7846433333333334
$ws.Range("H2").Value = 2.35393
$ws.Range("I2").Value = 0.8694484673945282
$ws.Range("J2").Value = 0.8694484673945282
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 173.637756
$ws.Range("N2").Value = 520.913268
$ws.Range("O2").Value = 0.9875738673498291
$ws.Range("P2").Value = 0.9875738673498291
$ws.Range("Q2").Value = 136.24370766036
$ws.Range("R2").Value = 1226.19336894324
$ws.Range("S2").Value = 0.858644585406196
$ws.Range("T2").Value = 0.858644585406196

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "F10"
$ws.Range("C3").Value = "F3"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7846433333333334
$ws.Range("H3").Value = 2.35393
$ws.Range("I3").Value = 0.8694484673945282
$ws.Range("J3").Value = 0.8694484673945282
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.184794333333333
$ws.Range("N3").Value = 6.554383
$ws.Range("O3").Value = 0.0124261326501708
$ws.Range("P3").Value = 0.0124261326501708
$ws.Range("Q3").Value = 1.714284308354444
$ws.Range("R3").Value = 15.42855877519
$ws.Range("S3").Value = 0.01080388198833211
$ws.Range("T3").Value = 0.01080388198833211

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "F10"
$ws.Range("C4").Value = "F3"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.09213
$ws.Range("H4").Value = 0.27639
$ws.Range("I4").Value = 0.1020875140310772
$ws.Range("J4").Value = 0.1020875140310772
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 173.637756
$ws.Range("N4").Value = 520.913268
$ws.Range("O4").Value = 0.9875738673498291
$ws.Range("P4").Value = 0.9875738673498291
$ws.Range("Q4").Value = 15.99724646028
$ws.Range("R4").Value = 143.97521814252
$ws.Range("S4").Value = 0.1008189610398009
$ws.Range("T4").Value = 0.1008189610398009

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "F10"
$ws.Range("C5").Value = "F3"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.09213
$ws.Range("H5").Value = 0.27639
$ws.Range("I5").Value = 0.1020875140310772
$ws.Range("J5").Value = 0.1020875140310772
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.184794333333333
$ws.Range("N5").Value = 6.554383
$ws.Range("O5").Value = 0.0124261326501708
$ws.Range("P5").Value = 0.0124261326501708
$ws.Range("Q5").Value = 0.20128510193
$ws.Range("R5").Value = 1.81156591737
$ws.Range("S5").Value = 0.001268552991276339
$ws.Range("T5").Value = 0.001268552991276339

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "F10"
$ws.Range("C6").Value = "F3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.02568766666666667
$ws.Range("H6").Value = 0.077063
$ws.Range("I6").Value = 0.02846401857439453
$ws.Range("J6").Value = 0.02846401857439454
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 173.637756
$ws.Range("N6").Value = 520.913268
$ws.Range("O6").Value = 0.9875738673498291
$ws.Range("P6").Value = 0.9875738673498291
$ws.Range("Q6").Value = 4.460348796876
$ws.Range("R6").Value = 40.143139171884
$ws.Range("S6").Value = 0.02811032090383218
$ws.Range("T6").Value = 0.02811032090383218

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "F10"
$ws.Range("C7").Value = "F3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.02568766666666667
$ws.Range("H7").Value = 0.077063
$ws.Range("I7").Value = 0.02846401857439453
$ws.Range("J7").Value = 0.02846401857439454
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.184794333333333
$ws.Range("N7").Value = 6.554383
$ws.Range("O7").Value = 0.0124261326501708
$ws.Range("P7").Value = 0.0124261326501708
$ws.Range("Q7").Value = 0.05612226856988888
$ws.Range("R7").Value = 0.505100417129
$ws.Range("S7").Value = 0.0003536976705623521
$ws.Range("T7").Value = 0.0003536976705623522
